$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: STAFF_ID changes 1120 -> 110; NO column (A12) gets a Text number format applied (value stays numeric 1) ---
$ws.Range("A12").NumberFormat = "@"
$ws.Range("B12").Value = 110

# --- Row 13: NO becomes text "2"; STAFF_ID changes 110 -> 111 ---
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "2"
$ws.Range("B13").Value = 111

# --- Row 14 (new staff member: Mark Mensah, Male) ---
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "3"
$ws.Range("B14").Value = 112
$ws.Range("C14").Value = "Mensah"
$ws.Range("D14").Value = "Mark "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0249626071"
$ws.Range("F14").Value = "Male"
$ws.Range("G14").Value = "yes"

# --- Row 15 (new staff member: Solace Dotsey, Female) ---
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "4"
$ws.Range("B15").Value = 113
$ws.Range("C15").Value = "Dotsey"
$ws.Range("D15").Value = "Solace"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0249626071"
$ws.Range("F15").Value = "Female"
$ws.Range("G15").Value = "yes"

# Match the recorded cursor position left behind by the edit.
$ws.Range("A12").Select()
